$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting already used for column A1 (bold, border, centered)
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill team record values for every data row (2-63)
$lastRow = 63
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 95
    $ws.Cells.Item($r, 31).Value = 67
    $ws.Cells.Item($r, 32).Value = 0
}
